$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between row 2 and row 3 for columns A, Q, R
$a2 = $ws.Range("A2").Value2
$a3 = $ws.Range("A3").Value2
$ws.Range("A2").Value = $a3
$ws.Range("A3").Value = $a2

$q2 = $ws.Range("Q2").Value2
$q3 = $ws.Range("Q3").Value2
$ws.Range("Q2").Value = $q3
$ws.Range("Q3").Value = $q2

$r2 = $ws.Range("R2").Value2
$r3 = $ws.Range("R3").Value2
$ws.Range("R2").Value = $r3
$ws.Range("R3").Value = $r2

# Swap values between row 5 and row 6 for columns A, Q, R, AC
$a5 = $ws.Range("A5").Value2
$a6 = $ws.Range("A6").Value2
$ws.Range("A5").Value = $a6
$ws.Range("A6").Value = $a5

$q5 = $ws.Range("Q5").Value2
$q6 = $ws.Range("Q6").Value2
$ws.Range("Q5").Value = $q6
$ws.Range("Q6").Value = $q5

$r5 = $ws.Range("R5").Value2
$r6 = $ws.Range("R6").Value2
$ws.Range("R5").Value = $r6
$ws.Range("R6").Value = $r5

$ac5 = $ws.Range("AC5").Value2
$ac6 = $ws.Range("AC6").Value2
$ws.Range("AC5").Value = $ac6
$ws.Range("AC6").Value = $ac5
